$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. New "Abstract Title" paragraph style (styleId "AbstractTitle"), based on
#    Normal, followed by the existing "Abstract" style.
# ---------------------------------------------------------------------------
$abstractTitle = $d.Styles.Add("Abstract Title", 1)
$abstractTitle.BaseStyle = $d.Styles("Normal")
$abstractTitle.NextParagraphStyle = $d.Styles("Abstract")
$abstractTitle.QuickStyle = $true

$abstractTitle.ParagraphFormat.KeepWithNext = $true
$abstractTitle.ParagraphFormat.KeepTogether = $true
$abstractTitle.ParagraphFormat.Alignment = 1
$abstractTitle.ParagraphFormat.SpaceBefore = 15
$abstractTitle.ParagraphFormat.SpaceAfter = 0

$abstractTitle.Font.Size = 10
$abstractTitle.Font.SizeBi = 10
$abstractTitle.Font.Bold = $true
$abstractTitle.Font.Color = 9067060

# ---------------------------------------------------------------------------
# 2. "Abstract" style: tighten the space that precedes it (before: 300 -> 100
#    twips, i.e. 15pt -> 5pt); the trailing space stays at 300 twips (15pt).
# ---------------------------------------------------------------------------
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5
$abstract.ParagraphFormat.SpaceAfter = 15

# ---------------------------------------------------------------------------
# 3. New "Footnote Block Text" paragraph style (styleId "FootnoteBlockText"),
#    based on "Footnote Text", mirroring the indentation/spacing of the
#    existing "Block Text" style.
# ---------------------------------------------------------------------------
$footnoteBlockText = $d.Styles.Add("Footnote Block Text", 1)
$footnoteBlockText.BaseStyle = $d.Styles("Footnote Text")
$footnoteBlockText.NextParagraphStyle = $d.Styles("Footnote Text")
$footnoteBlockText.Priority = 9
$footnoteBlockText.UnhideWhenUsed = $true
$footnoteBlockText.QuickStyle = $true

$footnoteBlockText.ParagraphFormat.SpaceBefore = 5
$footnoteBlockText.ParagraphFormat.SpaceAfter = 5
$footnoteBlockText.ParagraphFormat.FirstLineIndent = 0
$footnoteBlockText.ParagraphFormat.LeftIndent = 24
$footnoteBlockText.ParagraphFormat.RightIndent = 24

Write-Output "Abstract Title style: $($abstractTitle.NameLocal)"
Write-Output "Abstract style before-spacing: $($abstract.ParagraphFormat.SpaceBefore)"
Write-Output "Footnote Block Text style: $($footnoteBlockText.NameLocal)"
